$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update header row (row 1) values for columns B-E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 (CON) values for columns B-E
$ws.Range("B2").Value = 8.2959479899363604
$ws.Range("C2").Value = 5.6120728701218381
$ws.Range("D2").Value = 7.0690348294712946
$ws.Range("E2").Value = 7.7660648385154882

# Update row 3 (STR) values for columns B-E
$ws.Range("B3").Value = 7.0303179445172486
$ws.Range("C3").Value = 4.8852490163363234
$ws.Range("D3").Value = 5.8366287550394942
$ws.Range("E3").Value = 5.5878078263626962

# Update the selected range to match the new selection noted in the diff
$ws.Range("B1:E3").Select()
